# Auto-generated script applying the numeric/text updates from the commit diff
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 80
$ws.Range("F3").Value = 121
$ws.Range("F4").Value = 627
$ws.Range("F5").Value = 346
$ws.Range("F6").Value = 542
$ws.Range("F8").Value = 11248
$ws.Range("F12").Value = 2075
$ws.Range("F16").Value = 221
$ws.Range("F18").Value = 1170
$ws.Range("F19").Value = 147
$ws.Range("F20").Value = 234
$ws.Range("F21").Value = 719
$ws.Range("F23").Value = 244
$ws.Range("F24").Value = 2392
$ws.Range("F25").Value = 708
$ws.Range("F26").Value = 3380
$ws.Range("F27").Value = 1039
$ws.Range("F28").Value = 769
$ws.Range("F30").Value = 23
$ws.Range("F32").Value = 959
$ws.Range("F34").Value = 49
$ws.Range("F35").Value = 247
$ws.Range("F38").Value = 1625
$ws.Range("F39").Value = 4336
$ws.Range("F40").Value = 5431
$ws.Range("F42").Value = 108
$ws.Range("F43").Value = 22
$ws.Range("F44").Value = 144
$ws.Range("F45").Value = 244
$ws.Range("F46").Value = 57
$ws.Range("F47").Value = 20
$ws.Range("F48").Value = 4086

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 4130
$ws.Range("F5").Value = 83
$ws.Range("F7").Value = 36
$ws.Range("F11").Value = 522

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 747
$ws.Range("C4").Value = "杭州·GOGOGOODS谷子快跑（免费入场）"

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 747
$ws.Range("F5").Value = 80
$ws.Range("F6").Value = 627
$ws.Range("F7").Value = 346
$ws.Range("F8").Value = 542
$ws.Range("F9").Value = 11251
$ws.Range("F12").Value = 2075
$ws.Range("F15").Value = 221
$ws.Range("F17").Value = 1170
$ws.Range("F18").Value = 147
$ws.Range("F19").Value = 234
$ws.Range("F20").Value = 4130
$ws.Range("F22").Value = 719
$ws.Range("F24").Value = 244
$ws.Range("F25").Value = 708
$ws.Range("F26").Value = 1039
$ws.Range("F27").Value = 83
$ws.Range("F29").Value = 769
$ws.Range("F31").Value = 23
$ws.Range("F33").Value = 49
$ws.Range("F37").Value = 108
$ws.Range("F38").Value = 144
$ws.Range("F39").Value = 244
$ws.Range("F41").Value = 57
$ws.Range("F42").Value = 20
$ws.Range("F43").Value = 4086
$ws.Range("C4").Value = "杭州·GOGOGOODS谷子快跑（免费入场）"
